# Updated symbol list on Mon Dec 19 09:51:19 UTC 2022 with GitHub Actions
# Applies updated Price (column D) and Volume(1h) (column E) values to the
# cryptos worksheet. All these cells hold text (not numeric) values in the
# source workbook, so a leading apostrophe is used to force text storage
# and avoid Excel auto-converting numeric-looking strings to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "247.62"
    "D4"  = "5.472"
    "D5"  = "0.05695"
    "D6"  = "3.383"
    "D7"  = "0.8055"
    "D8"  = "1.039"
    "D9"  = "0.1468"
    "D10" = "0.07344"
    "D11" = "0.03167"
    "D12" = "0.02944"
    "D13" = "0.09292"
    "D14" = "0.001642"
    "D15" = "3.365"
    "D16" = "0.04700"
    "D17" = "0.0005871"
    "E17" = "16OneONE"
    "D18" = "0.006334"
    "D19" = "0.005045"
    "E19" = "18HotbitTokenHTB"
    "D20" = "0.001046"
    "D22" = "0.0003136"
    "D24" = "6.425"
    "D25" = "2.107"
    "D26" = "0.3289"
    "D40" = "0.04104"
    "D41" = "0.006942"
    "E41" = "40KickTokenKICKBestin24h"
    "D43" = "0.1044"
    "D44" = "0.008144"
    "D45" = "0.00005842"
    "D47" = "0.0005501"
    "E47" = "46ACDXExchangeACXTWorstin24h"
    "D48" = "0.6826"
    "D49" = "0.009585"
    "D50" = "0.00002100"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = "'" + $updates[$addr]
}
